# Auto-generated edit script applying numeric corrections to the
# "Phantom_Profits" workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2231.077
$ws.Range("I15").Value = 2231.077
$ws.Range("K15").Value = 6693.231000000001
$ws.Range("M15").Value = -6524.231000000001
$ws.Range("H21").Value = 29999
$ws.Range("I21").Value = 29999
$ws.Range("K21").Value = 29999
$ws.Range("M21").Value = -29531
$ws.Range("H23").Value = 29999
$ws.Range("I23").Value = 29999
$ws.Range("K23").Value = 29999
$ws.Range("M23").Value = -29765
$ws.Range("H40").Value = 1390.909
$ws.Range("I40").Value = 1385.5714
$ws.Range("K40").Value = 1385.5714
$ws.Range("M40").Value = -1210.5714
$ws.Range("H43").Value = 932.5
$ws.Range("I43").Value = 879.2
$ws.Range("K43").Value = 879.2
$ws.Range("M43").Value = -810.2
$ws.Range("H98").Value = 1093.375
$ws.Range("I98").Value = 633
$ws.Range("K98").Value = 633
$ws.Range("M98").Value = 865
$ws.Range("H112").Value = 4541.8184
$ws.Range("J112").Value = 4631
$ws.Range("L112").Value = 13893
$ws.Range("N112").Value = -16109
$ws.Range("H122").Value = 1093.375
$ws.Range("I122").Value = 633
$ws.Range("K122").Value = 1899
$ws.Range("M122").Value = 551
$ws.Range("H137").Value = 3334.76
$ws.Range("I137").Value = 3324.5217
$ws.Range("J137").Value = 3452.5
$ws.Range("K137").Value = 9973.5651
$ws.Range("L137").Value = 10357.5
$ws.Range("M137").Value = -7423.5651
$ws.Range("N137").Value = -15457.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2340.257
$ws.Range("I61").Value = 2198.6206
$ws.Range("K61").Value = 2198.6206
$ws.Range("M61").Value = -1986.6206
$ws.Range("H63").Value = 12468.117
$ws.Range("J63").Value = 15649.75
$ws.Range("L63").Value = 15649.75
$ws.Range("N63").Value = -17021.75
$ws.Range("H66").Value = 12468.117
$ws.Range("J66").Value = 15649.75
$ws.Range("L66").Value = 78248.75
$ws.Range("N66").Value = -85112.75
$ws.Range("H132").Value = 2943.125
$ws.Range("I132").Value = 2788.3928
$ws.Range("J132").Value = 4026.25
$ws.Range("K132").Value = 8365.178400000001
$ws.Range("L132").Value = 12078.75
$ws.Range("M132").Value = -5835.178400000001
$ws.Range("N132").Value = -17138.75
$ws.Range("H136").Value = 2340.257
$ws.Range("I136").Value = 2198.6206
$ws.Range("K136").Value = 6595.861800000001
$ws.Range("M136").Value = -4045.861800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1852.6
$ws.Range("I20").Value = 1927.875
$ws.Range("J20").Value = 1551.5
$ws.Range("K20").Value = 1927.875
$ws.Range("L20").Value = 1551.5
$ws.Range("M20").Value = -1680.875
$ws.Range("N20").Value = -2045.5
$ws.Range("I22").Value = 459.66666
$ws.Range("J22").Value = 70
$ws.Range("K22").Value = 459.66666
$ws.Range("L22").Value = 70
$ws.Range("M22").Value = -286.66666
$ws.Range("N22").Value = -416
$ws.Range("H107").Value = 2092.5
$ws.Range("I107").Value = 1704.2858
$ws.Range("K107").Value = 1704.2858
$ws.Range("M107").Value = 215.7141999999999
$ws.Range("H134").Value = 4592.5
$ws.Range("I134").Value = 4732
$ws.Range("K134").Value = 14196
$ws.Range("M134").Value = -11661
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 617.875
$ws.Range("I16").Value = 617.875
$ws.Range("K16").Value = 617.875
$ws.Range("M16").Value = -330.875
$ws.Range("H62").Value = 3204.818
$ws.Range("I62").Value = 2586.5715
$ws.Range("J62").Value = 4286.75
$ws.Range("K62").Value = 2586.5715
$ws.Range("L62").Value = 4286.75
$ws.Range("M62").Value = -1962.5715
$ws.Range("N62").Value = -5534.75
$ws.Range("H65").Value = 3204.818
$ws.Range("I65").Value = 2586.5715
$ws.Range("J65").Value = 4286.75
$ws.Range("K65").Value = 12932.8575
$ws.Range("L65").Value = 21433.75
$ws.Range("M65").Value = -9812.8575
$ws.Range("N65").Value = -27673.75
$ws.Range("H86").Value = 9498.666999999999
$ws.Range("I86").Value = 9664.333000000001
$ws.Range("J86").Value = 9333
$ws.Range("K86").Value = 9664.333000000001
$ws.Range("L86").Value = 9333
$ws.Range("M86").Value = -8541.333000000001
$ws.Range("N86").Value = -11579
$ws.Range("H89").Value = 9498.666999999999
$ws.Range("I89").Value = 9664.333000000001
$ws.Range("J89").Value = 9333
$ws.Range("K89").Value = 48321.665
$ws.Range("L89").Value = 46665
$ws.Range("M89").Value = -42705.665
$ws.Range("N89").Value = -57897
$ws.Range("H99").Value = 7310
$ws.Range("I99").Value = 5875
$ws.Range("J99").Value = 8745
$ws.Range("K99").Value = 5875
$ws.Range("L99").Value = 8745
$ws.Range("M99").Value = -4377
$ws.Range("N99").Value = -11741
$ws.Range("H107").Value = 345.9091
$ws.Range("I107").Value = 280.5
$ws.Range("K107").Value = 280.5
$ws.Range("M107").Value = 1639.5
$ws.Range("H113").Value = 617.875
$ws.Range("I113").Value = 617.875
$ws.Range("K113").Value = 617.875
$ws.Range("M113").Value = 1552.125
$ws.Range("H120").Value = 39975
$ws.Range("J120").Value = 39975
$ws.Range("L120").Value = 39975
$ws.Range("N120").Value = -47233
$ws.Range("H122").Value = 3870
$ws.Range("I122").Value = 4115.4
$ws.Range("K122").Value = 12346.2
$ws.Range("M122").Value = -9896.199999999999
$ws.Range("H126").Value = 7310
$ws.Range("I126").Value = 5875
$ws.Range("J126").Value = 8745
$ws.Range("K126").Value = 17625
$ws.Range("L126").Value = 26235
$ws.Range("M126").Value = -15155
$ws.Range("N126").Value = -31175

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 19002.5
$ws.Range("I76").Value = 15000
$ws.Range("J76").Value = 19803
$ws.Range("K76").Value = 45000
$ws.Range("L76").Value = 59409
$ws.Range("M76").Value = -44617
$ws.Range("N76").Value = -60175
$ws.Range("H79").Value = 19002.5
$ws.Range("I79").Value = 15000
$ws.Range("J79").Value = 19803
$ws.Range("K79").Value = 45000
$ws.Range("L79").Value = 59409
$ws.Range("M79").Value = -43674
$ws.Range("N79").Value = -62061
$ws.Range("H92").Value = 293.8
$ws.Range("I92").Value = 293.8
$ws.Range("K92").Value = 881.4000000000001
$ws.Range("M92").Value = 366.5999999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8779.6
$ws.Range("I70").Value = 8646.333000000001
$ws.Range("K70").Value = 8646.333000000001
$ws.Range("M70").Value = -8376.333000000001
$ws.Range("H73").Value = 8779.6
$ws.Range("I73").Value = 8646.333000000001
$ws.Range("K73").Value = 8646.333000000001
$ws.Range("M73").Value = -7710.333000000001
$ws.Range("H102").Value = 2265.4167
$ws.Range("I102").Value = 2107.818
$ws.Range("K102").Value = 2107.818
$ws.Range("M102").Value = -485.8180000000002
$ws.Range("H122").Value = 2252.85
$ws.Range("I122").Value = 2124.9285
$ws.Range("J122").Value = 2551.3333
$ws.Range("K122").Value = 6374.7855
$ws.Range("L122").Value = 7653.999899999999
$ws.Range("M122").Value = -3924.7855
$ws.Range("N122").Value = -12553.9999
$ws.Range("H132").Value = 4612
$ws.Range("I132").Value = 4249.0835
$ws.Range("K132").Value = 12747.2505
$ws.Range("M132").Value = -10217.2505
$ws.Range("H134").Value = 48000
$ws.Range("J134").Value = 48000
$ws.Range("L134").Value = 144000
$ws.Range("N134").Value = -149070

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3344.625
$ws.Range("I7").Value = 3344.625
$ws.Range("K7").Value = 3344.625
$ws.Range("M7").Value = -3232.625
$ws.Range("H126").Value = 3344.625
$ws.Range("I126").Value = 3344.625
$ws.Range("K126").Value = 10033.875
$ws.Range("M126").Value = -7563.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 29475
$ws.Range("J46").Value = 29475
$ws.Range("L46").Value = 29475
$ws.Range("N46").Value = -29937
$ws.Range("H122").Value = 3063.652
$ws.Range("I122").Value = 3171.842
$ws.Range("J122").Value = 2549.75
$ws.Range("K122").Value = 9515.526
$ws.Range("L122").Value = 7649.25
$ws.Range("M122").Value = -7065.526
$ws.Range("N122").Value = -12549.25
$ws.Range("H132").Value = 4209.0347
$ws.Range("I132").Value = 3761.3704
$ws.Range("K132").Value = 11284.1112
$ws.Range("M132").Value = -8754.111199999999
$ws.Range("H134").Value = 29475
$ws.Range("J134").Value = 29475
$ws.Range("L134").Value = 88425
$ws.Range("N134").Value = -93495

